$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 60-61 (pushes former rows 60-73 down to 62-75),
# mirroring the weekly-update pattern: a fresh week's pricing data for
# "Cebollín" gets inserted in date order, ahead of the already-present
# rows for later dates that were appended out of order.
$ws.Rows("60:61").Insert()

# New row 60 - "Primera" quality, week of 2022-10-11 (serial 44845)
$ws.Cells.Item(60, 1).Value = 11
$ws.Cells.Item(60, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(60, 3).Value = "Bíobío"
$ws.Cells.Item(60, 4).Value = 44845
$ws.Cells.Item(60, 5).Value = 8
$ws.Cells.Item(60, 6).Value = 100112037
$ws.Cells.Item(60, 7).Value = "Cebollín"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 250
$ws.Cells.Item(60, 11).Value = 2000
$ws.Cells.Item(60, 12).Value = 2000
$ws.Cells.Item(60, 13).Value = 2000
$ws.Cells.Item(60, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(60, 15).Value = "Región Metropolitana"
$ws.Cells.Item(60, 16).Value = 56
$ws.Cells.Item(60, 17).Value = 36
$ws.Cells.Item(60, 18).Value = "Hortaliza"

# New row 61 - "Segunda" quality, same week (serial 44845)
$ws.Cells.Item(61, 1).Value = 11
$ws.Cells.Item(61, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(61, 3).Value = "Bíobío"
$ws.Cells.Item(61, 4).Value = 44845
$ws.Cells.Item(61, 5).Value = 8
$ws.Cells.Item(61, 6).Value = 100112037
$ws.Cells.Item(61, 7).Value = "Cebollín"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Segunda"
$ws.Cells.Item(61, 10).Value = 150
$ws.Cells.Item(61, 11).Value = 1600
$ws.Cells.Item(61, 12).Value = 1600
$ws.Cells.Item(61, 13).Value = 1600
$ws.Cells.Item(61, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(61, 15).Value = "Región Metropolitana"
$ws.Cells.Item(61, 16).Value = 44
$ws.Cells.Item(61, 17).Value = 36
$ws.Cells.Item(61, 18).Value = "Hortaliza"
